$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# URL re-used across many photo cells
$url = "https://www.flickr.com/photos/aaronhatcher/45160698381/in/datetaken/"

# --- Row 1 (headers) ---
$ws.Range("A1").Value = "Team"
$ws.Range("B1").Value = "Player_Name"
$ws.Range("C1").Value = "Number"
$ws.Range("D1").Value = "Photo1"
$ws.Range("E1").Value = "Photo2"
$ws.Range("F1").Value = "Photo3"
$ws.Range("G1").Value = "Photo4"
$ws.Range("H1").Value = "Photo5"

# --- Row 2: Rebels / Smoosh ---
$ws.Range("A2").Value = "Rebels"
$ws.Range("B2").Value = "Smoosh"
$ws.Range("C2").Value = 10
$ws.Range("D2:E2").Value = $url

# --- Row 3: JV / Alex ---
$ws.Range("A3").Value = "JV"
$ws.Range("B3").Value = "Alex"
$ws.Range("C3").Value = 19
$ws.Range("D3:F3").Value = $url

# --- Row 4: JV / Tom ---
$ws.Range("A4").Value = "JV"
$ws.Range("B4").Value = "Tom"
$ws.Range("C4").Value = 20
$ws.Range("D4:E4").Value = $url

# --- Row 5: Varsity / Dick ---
$ws.Range("A5").Value = "Varsity"
$ws.Range("B5").Value = "Dick"
$ws.Range("C5").Value = 1
$ws.Range("D5:H5").Value = $url

# --- Row 6: Varsity / Harry ---
$ws.Range("A6").Value = "Varsity"
$ws.Range("B6").Value = "Harry"
$ws.Range("C6").Value = 2
$ws.Range("D6:I6").Value = $url

# --- Styling ---
# D2:E2 and D3:F3 keep the original named "Hyperlink" style (fontId 1 / xfId 1)
$ws.Range("D2:E2").Style = "Hyperlink"
$ws.Range("D3:F3").Style = "Hyperlink"

# D4:I6 (the newer photo links) get a direct underline + custom blue font
# color (fontId 2 / xfId 0, applyFont) instead of the named Hyperlink style.
# Build the format once on a scratch cell, then paste *formats only* onto
# the target ranges so every cell keeps its own text value.
$scratch = $ws.Range("Z100")
$scratch.Value = "x"
$scratch.Font.Underline = 2
$scratch.Font.Color = 12673797
$scratch.Copy()

$ws.Range("D4:E4").PasteSpecial(-4122)
$ws.Range("D5:H5").PasteSpecial(-4122)
$ws.Range("D6:I6").PasteSpecial(-4122)

$scratch.Clear()
$excel.CutCopyMode = $false

# --- Selection / active cell matches the saved workbook state ---
$ws.Range("I6").Select()
